# Updated cryptos list with refreshed prices / volumes and a couple of
# row re-orderings (ARBITRUM/LidoDAOToken swap; BabyDogeCoin inserted,
# NEARProtocol dropped off the bottom of the table).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new value. Column D holds numeric-looking text (prices with
# trailing zeros / thousand-dot formatting like "29.386.86"); it is
# written as Text so Excel does not silently coerce it to a Double and
# drop significant trailing/leading zeros.
$updates = [ordered]@{
    "D2" = '29.386.86'
    "E2" = '  -0.04%  '
    "D3" = '1.848.70'
    "D4" = '0.9989'
    "E4" = '  -0.06%  '
    "D5" = '240.32'
    "E5" = '  -0.08%  '
    "D6" = '0.6275'
    "E6" = '  -0.38%  '
    "D7" = '1.0000'
    "E7" = '  -0.03%  '
    "D8" = '0.07613'
    "E8" = '  -0.82%  '
    "D9" = '0.2903'
    "E9" = '  -1.20%  '
    "D10" = '24.74'
    "E10" = '  +1.05%  '
    "D11" = '0.07737'
    "E11" = '  -0.11%  '
    "D12" = '5.032'
    "E12" = '  +0.39%  '
    "D13" = '0.6794'
    "D14" = '0.00001055'
    "E14" = '  -3.60%  '
    "D15" = '83.10'
    "E15" = '  -0.57%  '
    "D16" = '6.158'
    "E16" = '  +0.30%  '
    "D17" = '29.411.96'
    "E17" = '  -0.03%  '
    "D18" = '227.65'
    "E18" = '  -0.81%  '
    "E19" = '  -0.83%  '
    "E20" = '  -0.06%  '
    "D21" = '7.481'
    "E21" = '  +0.57%  '
    "D22" = '0.9991'
    "E22" = '  -0.12%  '
    "D23" = '158.69'
    "E23" = '  +0.89%  '
    "E24" = '  -0.31%  '
    "D25" = '8.412'
    "E25" = '  +0.58%  '
    "D26" = '17.68'
    "E26" = '  +0.07%  '
    "D27" = '1.410'
    "E27" = '  +7.55%  '
    "D28" = '1.461'
    "E28" = '  -0.53%  '
    "D29" = '0.05603'
    "E29" = '  -1.37%  '
    "D30" = '4.112'
    "E30" = '  +0.03%  '
    "B32" = 'LidoDAOToken'
    "C32" = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
    "D32" = '1.837'
    "E32" = '  -0.76%  '
    "B33" = 'ARBITRUM'
    "C33" = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
    "D33" = '1.164'
    "E33" = '  +0.46%  '
    "D34" = '0.6977'
    "E34" = '  -1.77%  '
    "D35" = '2.589'
    "E35" = '  +0.16%  '
    "E36" = '  +0.24%  '
    "D37" = '1.229.35'
    "E37" = '  -0.15%  '
    "E38" = '  -1.74%  '
    "D39" = '6.386'
    "E39" = '  -1.28%  '
    "D40" = '0.9000'
    "E40" = '  -1.52%  '
    "D41" = '1.000'
    "E41" = '  +0.00%  '
    "D42" = '101.47'
    "E42" = '  +0.09%  '
    "D43" = '65.90'
    "E43" = '  -0.40%  '
    "D44" = '7.216'
    "E44" = '  +0.72%  '
    "B45" = 'BabyDogeCoin'
    "C45" = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
    "D45" = '0.00000000119'
    "E45" = '  -1.90%  '
    "B46" = 'TheSandbox'
    "C46" = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
    "D46" = '0.4001'
    "E46" = '  -0.29%  '
    "B47" = 'EnergySwap'
    "C47" = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
    "D47" = '9.005'
    "E47" = '  -0.58%  '
    "B48" = 'RenderToken'
    "C48" = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
    "D48" = '1.676'
    "E48" = '  -0.83%  '
    "B49" = 'Algorand'
    "C49" = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
    "D49" = '0.1138'
    "E49" = '  +1.19%  '
    "B50" = 'Cronos'
    "C50" = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
    "D50" = '0.05706'
    "E50" = '  -0.13%  '
    "B51" = 'Mantle'
    "C51" = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
    "D51" = '0.4627'
    "E51" = '  +0.03%  '
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    if ($ref[0] -eq "D") {
        $cell.NumberFormat = "@"
        $cell.Value = $updates[$ref]
        $cell.Style = "Normal"
    } else {
        $cell.Value = $updates[$ref]
    }
}
